$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")
$ws.Activate()

# Update the ModelParameterSheets value for TestScenario2 (row 3, column E):
# parenthesis-free parameter sheet list now includes a sheet name containing a comma.
$ws.Range("E3").Value = '"Global", "Aciclovir", "Sheet, with comma"'

# Move the active selection on the Scenarios sheet to E4
$ws.Range("E4").Select()
